$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 10.3
$ws.Range("E9").Value = 11.1
$ws.Range("F9").Value = 10.199999999999999
$ws.Range("H9").Value = 10.6
$ws.Range("I9").Value = 10.8
$ws.Range("F10").Value = 14.8
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = 14.1
$ws.Range("F11").Value = 7.2
$ws.Range("H11").Value = 8.4
$ws.Range("I11").Value = 8.6999999999999993
$ws.Range("D12").Value = 6.8
$ws.Range("E12").Value = 6.9
$ws.Range("F12").Value = 7
$ws.Range("H12").Value = 8.1999999999999993
$ws.Range("I12").Value = 8.4
$ws.Range("F13").Value = 8.5
$ws.Range("H13").Value = 9.4
$ws.Range("I13").Value = 9
$ws.Range("F14").Value = 5.9
$ws.Range("H14").Value = 7.3
$ws.Range("I14").Value = 7.9
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 7.2
$ws.Range("F15").Value = 9.4
$ws.Range("H15").Value = 9.6999999999999993
$ws.Range("I15").Value = 9.6
$ws.Range("F16").Value = 11.1
$ws.Range("H16").Value = 10.8
$ws.Range("I16").Value = 8.9
$ws.Range("F17").Value = 8.1999999999999993
$ws.Range("H17").Value = 8.9
$ws.Range("I17").Value = 10
$ws.Range("D18").Value = 11.1
$ws.Range("E18").Value = 11.2
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = 10.3
$ws.Range("I18").Value = 10.199999999999999
$ws.Range("F19").Value = 14.7
$ws.Range("H19").Value = 14.7
$ws.Range("I19").Value = 16.3
$ws.Range("F20").Value = 7.6
$ws.Range("H20").Value = 8.1
$ws.Range("I20").Value = 7.4
$ws.Range("D21").Value = 5.4
$ws.Range("E21").Value = 5.5
$ws.Range("F21").Value = 6.6
$ws.Range("H21").Value = 6.6
$ws.Range("I21").Value = 6.8
$ws.Range("F22").Value = 8.1
$ws.Range("H22").Value = 7.8
$ws.Range("I22").Value = 8.9
$ws.Range("F23").Value = 5.7
$ws.Range("H23").Value = 5.8
$ws.Range("I23").Value = 5.5
$ws.Range("D24").Value = 5.0999999999999996
$ws.Range("E24").Value = 5.4
$ws.Range("F24").Value = 4.9000000000000004
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 4.0999999999999996
$ws.Range("F25").Value = 5.7
$ws.Range("H25").Value = 5.9
$ws.Range("I25").Value = 4.0999999999999996
$ws.Range("F26").Value = 4.2
$ws.Range("H26").Value = 4.4000000000000004
$ws.Range("I26").Value = 4.0999999999999996
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = 12.7
$ws.Range("F27").Value = 10.3
$ws.Range("H27").Value = 9.9
$ws.Range("I27").Value = 9.9
$ws.Range("F28").Value = 11.3
$ws.Range("H28").Value = 10.5
$ws.Range("I28").Value = 9.1
$ws.Range("F29").Value = 9.6
$ws.Range("H29").Value = 9.5
$ws.Range("I29").Value = 10.5
$ws.Range("D30").Value = 8.6999999999999993
$ws.Range("E30").Value = 8.6
$ws.Range("F30").Value = 10.199999999999999
$ws.Range("H30").Value = 9.6
$ws.Range("I30").Value = 8.6999999999999993
$ws.Range("F31").Value = 9.6
$ws.Range("H31").Value = 10.6
$ws.Range("I31").Value = 9.6999999999999993
$ws.Range("F32").Value = 10.7
$ws.Range("H32").Value = 8.6999999999999993
$ws.Range("I32").Value = 7.8
$ws.Range("D33").Value = "…"
$ws.Range("E33").Value = "…"
$ws.Range("F33").Value = "…"
$ws.Range("D34").Value = "…"
$ws.Range("E34").Value = "…"
$ws.Range("F34").Value = "…"
$ws.Range("D35").Value = "…"
$ws.Range("E35").Value = "…"
$ws.Range("F35").Value = "…"

$ws.Range("C1").Select()

Write-Host "Edit complete"
